# Test Summary Report_V1.0_ClipShot.docx edit
# Commit message: "conversione documenti in pdf"
#
# The underlying XML diff shows the author re-typed/edited several spans of
# text (which merges the previously split "ClipShot" runs with their
# neighbours) and changed two numbers:
#   - "27" -> "21"  (pianificati: casi di test)
#   - "24" -> "18"  (Failed: numero di test case)
# plus the Word-managed "_GoBack" bookmark (marks the last edit point) moved
# from the end of the document to right after the new "18".

$d = $word.ActiveDocument

# --- Paragraph 1: merge the first (unflagged) "ClipShot" run into its
#     surrounding text (the run that stays wrapped in proofErr spellStart/
#     spellEnd a little further on is left alone). -------------------------
$d.Content.Find.Execute(
    "Nel seguente documento viene mostrato il resoconto delle attività di testing ed i relativi risultati finali fornendo una valutazione relativamente all'esecuzione dei test case specificati nel documento Test Case Specifications_V1.0_ClipShot. Effettuiamo dunque, una prima analisi dei difetti riscontrati nell'applicazione, ovvero un insieme di valutazioni rispetto ai fallimenti o gli errori riscontrati nei test, che necessitano di ulteriori approfondimenti. È necessario fornire un report significativo delle attività di testing svolte in precedenza in modo da poter mettere a conoscenza chi usufruirà di ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nel seguente documento viene mostrato il resoconto delle attività di testing ed i relativi risultati finali fornendo una valutazione relativamente all'esecuzione dei test case specificati nel documento Test Case Specifications_V1.0_ClipShot. Effettuiamo dunque, una prima analisi dei difetti riscontrati nell'applicazione, ovvero un insieme di valutazioni rispetto ai fallimenti o gli errori riscontrati nei test, che necessitano di ulteriori approfondimenti. È necessario fornire un report significativo delle attività di testing svolte in precedenza in modo da poter mettere a conoscenza chi usufruirà di ",
    2) | Out-Null

# --- Same paragraph: merge " Report_V1.0_" + "ClipShot" + the trailing
#     sentence into a single run. ------------------------------------------
$d.Content.Find.Execute(
    " Report_V1.0_ClipShot. Nel seguente documento, ci serviamo di tali esiti positivi o negativi per poter giungere a soluzioni ai limiti riscontrati nei test. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Report_V1.0_ClipShot. Nel seguente documento, ci serviamo di tali esiti positivi o negativi per poter giungere a soluzioni ai limiti riscontrati nei test. ",
    2) | Out-Null

# --- "Test Plan_V1.0_" + "ClipShot" + ", attraverso la tecnica del " -------
$d.Content.Find.Execute(
    "All'interno del documento Test Plan_V1.0_ClipShot, attraverso la tecnica del ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "All'interno del documento Test Plan_V1.0_ClipShot, attraverso la tecnica del ",
    2) | Out-Null

# --- "Test Case Specification_V1.0_" + "ClipShot" + the following sentence
#     up to (not including) the "27" number run. ---------------------------
$d.Content.Find.Execute(
    ", sono state definite le varie combinazioni per i possibili input all'interno del sistema. Successivamente, nel documento Test Case Specification_V1.0_ClipShot sono stati specificati in dettaglio i vari test case con il relativo comportamento atteso. Avevamo pianificato esattamente ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", sono state definite le varie combinazioni per i possibili input all'interno del sistema. Successivamente, nel documento Test Case Specification_V1.0_ClipShot sono stati specificati in dettaglio i vari test case con il relativo comportamento atteso. Avevamo pianificato esattamente ",
    2) | Out-Null

# Number of planned test cases: 27 -> 21
$planned = $d.Content
$planned.Find.Execute("27")
$planned.Text = "21"

# Number of Failed test cases: 24 -> 18
$failedCount = $d.Content
$failedCount.Find.Execute("24")
$failedCount.Text = "18"

# Re-create the "_GoBack" bookmark right after the new "18" (Word always
# tracks the last edited spot with this hidden bookmark; here it ends up
# between "18" and the following ". Soltanto ...").
$goBackSpot = $d.Content
$goBackSpot.Find.Execute("18. Soltanto")
$goBackRange = $d.Range($goBackSpot.Start + 2, $goBackSpot.Start + 2)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# --- Final paragraph: merge the "TC_REGISTRAZIONE_02, " run with the
#     trailing "TC_AggiungiPost_02, TC_RicercaUtente_02." run (the old
#     trailing "_GoBack" bookmark that used to sit here is gone now). ------
$d.Content.Find.Execute(
    "I test case in questione sono: TC_REGISTRAZIONE_02, TC_AggiungiPost_02, TC_RicercaUtente_02.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I test case in questione sono: TC_REGISTRAZIONE_02, TC_AggiungiPost_02, TC_RicercaUtente_02.",
    2) | Out-Null
